# Swap the two sheet tab names: "test 640" <-> "test 7418" while keeping
# each tab's underlying data (sheet1.xml / sheet2.xml) in place.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "__tmp_swap__"
$ws2.Name = "test 640"
$ws1.Name = "test 7418"

# $ws1 is the sheet now named "test 7418" (first tab) - update its data row.
# Writing text that looks like a date/number through COM auto-formats the
# cell (adds a style + may coerce type), so force Text format, assign the
# value, then clear the formatting that step introduced to land back on the
# default (unstyled) cell - matching how the source file stored these as
# plain shared-string text.
$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "2020-12-05"
$ws1.Range("A2").ClearFormats()

$ws1.Range("B2").Value = 25150.0

$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C2").Value = "2021-12-05"
$ws1.Range("C2").ClearFormats()

$ws1.Range("D2").NumberFormat = "@"
$ws1.Range("D2").Value = "1"
$ws1.Range("D2").ClearFormats()

# $ws2 is the sheet now named "test 640" (second tab) - its single data row
# is removed entirely, leaving only the header row.
$ws2.Rows(2).Delete()
